# "added 4wk low sales check"
# Refresh the forecast figures on the "Forecast Comparison" sheet (the
# forecast model was rerun after adding a 4-week rolling low-sales check,
# which pulled several weeks' forecasts down) and roll the new totals up
# into the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Forecast Comparison sheet — per-week MyForecast / Inventory Coverage /
# Stockout Risk / Reorder Urgency / Seasonality Index
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Week W10 (row 2)
$ws.Range("D2").Value = 91
$ws.Range("H2").Value = 11.83
$ws.Range("L2").Value = 1.2

# Week W11 (row 3)
$ws.Range("D3").Value = 91
$ws.Range("H3").Value = 10.81
$ws.Range("L3").Value = 1.01

# Week W12 (row 4)
$ws.Range("D4").Value = 91
$ws.Range("H4").Value = 9.87
$ws.Range("L4").Value = 1.12

# Week W13 (row 5)
$ws.Range("D5").Value = 90
$ws.Range("H5").Value = 8.99
$ws.Range("L5").Value = 1.05

# Week W14 (row 6)
$ws.Range("D6").Value = 88
$ws.Range("H6").Value = 8.17
$ws.Range("L6").Value = 1.08

# Week W15 (row 7)
$ws.Range("D7").Value = 85
$ws.Range("H7").Value = 7.43
$ws.Range("L7").Value = 1.1

# Week W16 (row 8)
$ws.Range("D8").Value = 83
$ws.Range("H8").Value = 6.53
$ws.Range("L8").Value = 1.09

# Week W17 (row 9)
$ws.Range("D9").Value = 84
$ws.Range("H9").Value = 5.48
$ws.Range("L9").Value = 1.08

# Week W18 (row 10)
$ws.Range("D10").Value = 86
$ws.Range("H10").Value = 4.39
$ws.Range("L10").Value = 0.85

# Week W19 (row 11) — inventory coverage recovers above the stockout
# threshold, so risk/urgency drop from High/Urgent to Low/Normal
$ws.Range("D11").Value = 86
$ws.Range("H11").Value = 3.36
$ws.Range("I11").Value = "Low"
$ws.Range("J11").Value = "Normal"
$ws.Range("L11").Value = 1.12

# Week W20 (row 12)
$ws.Range("D12").Value = 85
$ws.Range("H12").Value = 2.41
$ws.Range("I12").Value = "Low"
$ws.Range("J12").Value = "Normal"
$ws.Range("L12").Value = 1.09

# Week W21 (row 13)
$ws.Range("D13").Value = 82
$ws.Range("H13").Value = 1.45
$ws.Range("I13").Value = "Low"
$ws.Range("J13").Value = "Normal"
$ws.Range("L13").Value = 1.1

# Week W22 (row 14) — Stockout Risk/Reorder Urgency stay High/Urgent
$ws.Range("D14").Value = 81
$ws.Range("H14").Value = 0.46
$ws.Range("L14").Value = 0.94

# Week W23 (row 15)
$ws.Range("D15").Value = 81
$ws.Range("L15").Value = 1.07

# Week W24 (row 16)
$ws.Range("D16").Value = 81
$ws.Range("L16").Value = 1.01

# Week W25 (row 17)
$ws.Range("D17").Value = 81
$ws.Range("L17").Value = 0.93

# ---------------------------------------------------------------------
# Summary sheet — roll up the revised forecast totals. Column B holds
# text values (e.g. "1374"), so briefly force the Text number format
# before assigning the numeric-looking strings to keep them stored as
# text rather than being auto-converted to numbers, then clear the
# formatting back off again so the cell style is left untouched (same
# General-format, default-style text cell as before the edit).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B9").NumberFormat = "@"
$summary.Range("B9").Value = "1374"   # Total Forecast (16 Weeks)
$summary.Range("B9").ClearFormats()

$summary.Range("B10").NumberFormat = "@"
$summary.Range("B10").Value = "706"   # Total Forecast (8 Weeks)
$summary.Range("B10").ClearFormats()

$summary.Range("B11").NumberFormat = "@"
$summary.Range("B11").Value = "365"   # Total Forecast (4 Weeks)
$summary.Range("B11").ClearFormats()

$summary.Range("B12").NumberFormat = "@"
$summary.Range("B12").Value = "92"    # Max Forecast
$summary.Range("B12").ClearFormats()

$summary.Range("B14").NumberFormat = "@"
$summary.Range("B14").Value = "82"    # Min Forecast
$summary.Range("B14").ClearFormats()
